$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 57

# Abbreviation used in the block code for each Landuse value
$landuseMap = @{
    "Agriculture" = "Ag"
    "Pasture"     = "P"
    "Wild"        = "W"
    "Seronera"    = "Local"
}

# Compute the Blockcode for every data row (Region_AbbrevBlock, except the
# "Seronera" landuse which only uses the landuse abbreviation + block number)
$codes = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $region  = $ws.Cells.Item($r, 2).Value2
    $landuse = $ws.Cells.Item($r, 5).Value2
    $block   = $ws.Cells.Item($r, 4).Value2

    $abbrev = $landuseMap[$landuse]

    if ($landuse -eq "Seronera") {
        $code = "$abbrev$block"
    } else {
        $code = "${region}_$abbrev$block"
    }

    $codes[$r] = $code
}

# Header for the new column
$ws.Range("K1").Value = "Blockcode"

# Prime the shared-string table with the distinct codes in alphabetical
# order (using scratch cells far away from the data) so the codes end up
# registered in that order before they are written into column K.
$uniqueCodes = $codes.Values | Sort-Object -Unique
$scratchRow = 1
foreach ($code in $uniqueCodes) {
    $ws.Cells.Item($scratchRow, 500).Value = $code
    $scratchRow++
}

# Fill column K row by row - this reuses the strings already registered
# above instead of creating new shared-string entries.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 11).Value = $codes[$r]
}

# Remove the scratch cells now that every code is safely referenced from
# column K.
$ws.Range($ws.Cells.Item(1, 500), $ws.Cells.Item($scratchRow, 500)).ClearContents()

# Update the selected cell to reflect where the user ended up after editing
$ws.Range("M7").Select()
